$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column E (部门ID) top-to-bottom first
$ws.Range("E1").Value = "部门ID"
$ws.Range("E2").Value = 4
$ws.Range("E3").Value = 4
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 5

# Then fill column F (部门名称) top-to-bottom
$ws.Range("F1").Value = "部门名称"
$ws.Range("F2").Value = "销售一部"
$ws.Range("F3").Value = "销售一部"
$ws.Range("F4").Value = "销售二部"
$ws.Range("F5").Value = "销售二部"

# Then fill column G (岗位类型) top-to-bottom
$ws.Range("G1").Value = "岗位类型"
$ws.Range("G2").Value = "主管"
$ws.Range("G3").Value = "员工"
$ws.Range("G4").Value = "员工"
$ws.Range("G5").Value = "员工"

# Move selection to match the committed view state
$ws.Range("E5").Select()
